# Auto-generated Excel COM-interop script to apply Tonberry_Profits value updates
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 288
$ws.Range("I28").Value = 110.125
$ws.Range("K28").Value = 110.125
$ws.Range("M28").Value = 374.875
$ws.Range("H121").Value = 1212.6666
$ws.Range("J121").Value = 1744
$ws.Range("L121").Value = 5232
$ws.Range("N121").Value = -8726
$ws.Range("H132").Value = 6537324
$ws.Range("I132").Value = 7408678
$ws.Range("J132").Value = 2167
$ws.Range("K132").Value = 22226034
$ws.Range("L132").Value = 6501
$ws.Range("M132").Value = -22223504
$ws.Range("N132").Value = -11561
$ws.Range("H138").Value = 2070.6165
$ws.Range("I138").Value = 1803.6046
$ws.Range("J138").Value = 2453.3333
$ws.Range("K138").Value = 5410.8138
$ws.Range("L138").Value = 7359.999899999999
$ws.Range("M138").Value = -270.8137999999999
$ws.Range("N138").Value = -17639.9999
$ws.Range("H141").Value = 850199.3
$ws.Range("I141").Value = 934104.0600000001
$ws.Range("K141").Value = 2802312.18
$ws.Range("M141").Value = -2797132.18

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 3333
$ws.Range("I39").Value = 3333
$ws.Range("K39").Value = 3333
$ws.Range("M39").Value = -2813
$ws.Range("H45").Value = 1908.3334
$ws.Range("I45").Value = 1919
$ws.Range("K45").Value = 1919
$ws.Range("M45").Value = -1542
$ws.Range("H61").Value = 71432104
$ws.Range("I61").Value = 55558110
$ws.Range("K61").Value = 55558110
$ws.Range("M61").Value = -55557898
$ws.Range("H74").Value = 2153.0527
$ws.Range("I74").Value = 1833.8
$ws.Range("K74").Value = 1833.8
$ws.Range("M74").Value = -959.8
$ws.Range("H77").Value = 2153.0527
$ws.Range("I77").Value = 1833.8
$ws.Range("K77").Value = 9169
$ws.Range("M77").Value = -4801
$ws.Range("H132").Value = 1418.1
$ws.Range("I132").Value = 1035.2122
$ws.Range("J132").Value = 2161.353
$ws.Range("K132").Value = 3105.6366
$ws.Range("L132").Value = 6484.059
$ws.Range("M132").Value = -575.6365999999998
$ws.Range("N132").Value = -11544.059
$ws.Range("H136").Value = 71432104
$ws.Range("I136").Value = 55558110
$ws.Range("K136").Value = 166674330
$ws.Range("M136").Value = -166671780

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1092

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3573509
$ws.Range("I31").Value = 5496107
$ws.Range("J31").Value = 2970.4285
$ws.Range("K31").Value = 5496107
$ws.Range("L31").Value = 2970.4285
$ws.Range("M31").Value = -5495812
$ws.Range("N31").Value = -3560.4285
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N33").ClearContents()
$ws.Range("H34").Value = 3573509
$ws.Range("I34").Value = 5496107
$ws.Range("J34").Value = 2970.4285
$ws.Range("K34").Value = 5496107
$ws.Range("L34").Value = 2970.4285
$ws.Range("M34").Value = -5495905
$ws.Range("N34").Value = -3374.4285
$ws.Range("H35").Value = 362.5
$ws.Range("I35").Value = 362.5
$ws.Range("K35").Value = 362.5
$ws.Range("M35").Value = -68.5
$ws.Range("H58").Value = 5437182.5
$ws.Range("I58").Value = 10870440
$ws.Range("J58").Value = 3924.75
$ws.Range("K58").Value = 10870440
$ws.Range("L58").Value = 3924.75
$ws.Range("M58").Value = -10870237
$ws.Range("N58").Value = -4330.75
$ws.Range("H99").Value = 1662.4445
$ws.Range("I99").Value = 1620.25
$ws.Range("K99").Value = 1620.25
$ws.Range("M99").Value = -122.25
$ws.Range("H105").Value = 1369.5385
$ws.Range("I105").Value = 1391.0834
$ws.Range("K105").Value = 1391.0834
$ws.Range("M105").Value = 355.9166
$ws.Range("H126").Value = 1662.4445
$ws.Range("I126").Value = 1620.25
$ws.Range("K126").Value = 4860.75
$ws.Range("M126").Value = -2390.75
$ws.Range("H132").Value = 1564.8223
$ws.Range("I132").Value = 1148.0769
$ws.Range("K132").Value = 3444.2307
$ws.Range("M132").Value = -914.2307000000001
$ws.Range("H134").Value = 1813.9535
$ws.Range("I134").Value = 1613.1621
$ws.Range("K134").Value = 4839.4863
$ws.Range("M134").Value = -2304.4863
$ws.Range("H136").Value = 5437182.5
$ws.Range("I136").Value = 10870440
$ws.Range("J136").Value = 3924.75
$ws.Range("K136").Value = 32611320
$ws.Range("L136").Value = 11774.25
$ws.Range("M136").Value = -32608770
$ws.Range("N136").Value = -16874.25
$ws.Range("H138").Value = 84442
$ws.Range("J138").Value = 84442
$ws.Range("L138").Value = 84442
$ws.Range("N138").Value = -94722

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 7664.9287
$ws.Range("I56").Value = 7664.9287
$ws.Range("K56").Value = 7664.9287
$ws.Range("M56").Value = -7134.9287
$ws.Range("H113").Value = 7532.6665
$ws.Range("I113").Value = 50500.5
$ws.Range("J113").Value = 922.2308
$ws.Range("K113").Value = 151501.5
$ws.Range("L113").Value = 2766.6924
$ws.Range("M113").Value = -149331.5
$ws.Range("N113").Value = -7106.6924
$ws.Range("H122").Value = 878.6829
$ws.Range("I122").Value = 529.38464
$ws.Range("J122").Value = 1040.8572
$ws.Range("K122").Value = 4764.46176
$ws.Range("L122").Value = 9367.7148
$ws.Range("M122").Value = -2314.46176
$ws.Range("N122").Value = -14267.7148
$ws.Range("H131").Value = 5272592
$ws.Range("I131").Value = 62500532
$ws.Range("J131").Value = 10252.713
$ws.Range("K131").Value = 187501596
$ws.Range("L131").Value = 30758.139
$ws.Range("M131").Value = -187496556
$ws.Range("N131").Value = -40838.139

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 33339990
$ws.Range("I10").Value = 33339990
$ws.Range("K10").Value = 33339990
$ws.Range("M10").Value = -33339821
$ws.Range("H19").Value = 44668.332
$ws.Range("J19").Value = 38503
$ws.Range("L19").Value = 38503
$ws.Range("N19").Value = -39079
$ws.Range("H40").Value = 57018
$ws.Range("J40").Value = 57018
$ws.Range("L40").Value = 57018
$ws.Range("N40").Value = -57320
$ws.Range("H97").Value = 1277.8
$ws.Range("I97").Value = 1089.7693
$ws.Range("J97").Value = 2500
$ws.Range("K97").Value = 1089.7693
$ws.Range("L97").Value = 2500
$ws.Range("M97").Value = -593.7692999999999
$ws.Range("N97").Value = -3492
$ws.Range("H132").Value = 1834262.6
$ws.Range("I132").Value = 2960632.5
$ws.Range("J132").Value = 3911.625
$ws.Range("K132").Value = 8881897.5
$ws.Range("L132").Value = 11734.875
$ws.Range("M132").Value = -8879367.5
$ws.Range("N132").Value = -16794.875
$ws.Range("H136").Value = 8997.781999999999
$ws.Range("J136").Value = 8997.781999999999
$ws.Range("L136").Value = 26993.346
$ws.Range("N136").Value = -32093.346

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3186.9333
$ws.Range("I7").Value = 2607.7144
$ws.Range("K7").Value = 2607.7144
$ws.Range("M7").Value = -2495.7144
$ws.Range("H126").Value = 3186.9333
$ws.Range("I126").Value = 2607.7144
$ws.Range("K126").Value = 7823.1432
$ws.Range("M126").Value = -5353.1432
$ws.Range("H136").Value = 3063.138
$ws.Range("I136").Value = 1963.2609
$ws.Range("K136").Value = 5889.7827
$ws.Range("M136").Value = -3339.7827

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1721
$ws.Range("I132").Value = 1263.45
$ws.Range("K132").Value = 3790.35
$ws.Range("M132").Value = -1260.35
$ws.Range("H136").Value = 10895180
$ws.Range("J136").Value = 2120.6667
$ws.Range("L136").Value = 6362.000100000001
$ws.Range("N136").Value = -11462.0001
